# Strip stray leading/trailing spaces and surrounding quotation marks from
# the header/unit label cells in rows 1-3 (and the divider in row 4).
# Cells whose text was only spaces/quotes become empty and are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - main headers
$ws.Range("A1").Value = "TIME"
$ws.Range("B1").Value = "VO2"
$ws.Range("C1").Value = "VO2/kg"
$ws.Range("D1").Value = "METS"
$ws.Range("E1").Value = "VCO2"
$ws.Range("F1").Value = "VE"
$ws.Range("G1").Value = "RER"
$ws.Range("H1").Value = "RR"
$ws.Range("I1").Value = "Vt"
$ws.Range("J1").Value = "FEO2"
$ws.Range("K1").Value = "FECO2"
$ws.Range("L1").Value = "HR"
$ws.Range("M1").Value = "VE/"
$ws.Range("N1").Value = "VE/"
$ws.Range("O1").Value = "PetCO2"
$ws.Range("P1").Value = "Power"
$ws.Range("Q1").Value = "Lactate"

# Row 2 - secondary units (blank-only cells get cleared)
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "STPD"
$ws.Range("C2").Value = "STPD"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "STPD"
$ws.Range("F2").Value = "BTPS"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "BTPS"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "VO2"
$ws.Range("N2").Value = "VCO2"
$ws.Range("O2").Value = ""

# Row 3 - units
$ws.Range("A3").Value = "min"
$ws.Range("B3").Value = "L/min"
$ws.Range("C3").Value = "ml/kg/m"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "L/min"
$ws.Range("F3").Value = "L/min"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = "BPM"
$ws.Range("I3").Value = "L"
$ws.Range("J3").Value = "%"
$ws.Range("K3").Value = "%"
$ws.Range("L3").Value = "bpm"
$ws.Range("M3").Value = "BT/ST"
$ws.Range("N3").Value = "BT/ST"
$ws.Range("O3").Value = "mmHg"
$ws.Range("P3").Value = "W"
$ws.Range("Q3").Value = "mmol/L"

# Row 4 - divider
$ws.Range("A4").Value = "----------"
